# Insert two new highlighted paragraphs ("Say few words about .rodata..." and
# "Add the link to C++ standard...") right after the two blank paragraphs that
# follow "Slide 8." and immediately before the "Tell about ldd." paragraph.

$d = $word.ActiveDocument

# Locate the anchor paragraph ("Tell about ldd. In article on habr.") by text
# search so the script is resilient to the exact paragraph index.
$anchorRange = $d.Content
[void]$anchorRange.Find.Execute("Tell about")
$anchorPara = $anchorRange.Paragraphs(1)

# Make room for the new content: insert a fresh empty paragraph right before
# the anchor paragraph, then fill that empty paragraph with the exact OOXML
# for both new paragraphs (so run/highlight/proofErr structure matches
# precisely what Word itself would produce).
$insertionPoint = $anchorPara.Range
$insertionPoint.Collapse(1)
$insertionPoint.InsertParagraphBefore()

$newParaRange = $anchorRange.Paragraphs(1).Range

$xmlPayload = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>Say few words about .</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>rodata</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> section. UB when you try to </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>const_cast</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve"> it.</w:t></w:r></w:p><w:p><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>Add the link to C++ standard about one definition rule.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

[void]$newParaRange.InsertXML($xmlPayload)
